$d = $word.ActiveDocument
$vt = [char]0x0B

function Find-LabelEnd($range, $labelText) {
    $f = $range.Duplicate
    $f.Find.Execute($labelText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $f.End
}
function Find-LabelStart($range, $labelText) {
    $f = $range.Duplicate
    $f.Find.Execute($labelText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $f.Start
}
function Replace-RangeText($range, $oldText, $newText) {
    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# ------------------------------------------------------------------
# The document's content blocks rotate cyclically through the fixed
# sequence of headings:
#   Objetivos -> Docente(s) -> Programa resumido -> Programa ->
#   Avaliacao(Metodo) -> Avaliacao(Criterio) -> Avaliacao(Norma) ->
#   Bibliografia -> (back to Objetivos)
# i.e. each heading's body content is replaced by the content that,
# before the edit, belonged to the *next* heading in the list.
# ------------------------------------------------------------------

# ===================== CAPTURE PHASE (read only) ===================

# Paragraph 6  = Objetivos body
$t6 = $d.Paragraphs.Item(6).Range.Text
$t6 = $t6.Substring(0, $t6.Length - 1)

# Paragraph 8  = Docente(s) list body
$t8 = $d.Paragraphs.Item(8).Range.Text
$t8 = $t8.Substring(0, $t8.Length - 1)

# Paragraph 10 = Programa resumido body
$t10 = $d.Paragraphs.Item(10).Range.Text
$t10 = $t10.Substring(0, $t10.Length - 1)

# Paragraph 12 = Programa body
$t12 = $d.Paragraphs.Item(12).Range.Text
$t12 = $t12.Substring(0, $t12.Length - 1)

# Paragraph 14 = Avaliacao list (Metodo: / Criterio: / Norma de recuperacao:)
$p14 = $d.Paragraphs.Item(14).Range
$metodoLabelEnd   = Find-LabelEnd   $p14 "Método: "
$criterioStart    = Find-LabelStart $p14 "Critério: "
$criterioLabelEnd = Find-LabelEnd   $p14 "Critério: "
$normaStart       = Find-LabelStart $p14 "Norma de recuperação: "
$normaLabelEnd    = Find-LabelEnd   $p14 "Norma de recuperação: "
$p14End           = $p14.End

$oldMetodoContent   = $d.Range($metodoLabelEnd, $criterioStart).Text
$oldCriterioContent = $d.Range($criterioLabelEnd, $normaStart).Text
$oldNormaContentRaw = $d.Range($normaLabelEnd, $p14End).Text
$oldNormaContent    = $oldNormaContentRaw.Substring(0, $oldNormaContentRaw.Length - 1)

# Paragraph 16 = Bibliografia body
$t16 = $d.Paragraphs.Item(16).Range.Text
$t16 = $t16.Substring(0, $t16.Length - 1)

# ===================== WRITE PHASE (mutations) ======================

# Objetivos body <- old Programa resumido body
Replace-RangeText $d.Paragraphs.Item(6).Range $t6 $t10

# Docente(s) list body <- old Objetivos body
Replace-RangeText $d.Paragraphs.Item(8).Range $t8 $t6

# Programa resumido body <- old Programa body
Replace-RangeText $d.Paragraphs.Item(10).Range $t10 $t12

# Programa body <- old "Metodo:" content (strip its trailing manual break,
# since here it becomes a standalone paragraph with no following label)
$newProgramaContent = $oldMetodoContent
if ($newProgramaContent.Length -gt 0 -and [int][char]$newProgramaContent[$newProgramaContent.Length - 1] -eq 11) {
    $newProgramaContent = $newProgramaContent.Substring(0, $newProgramaContent.Length - 1)
}
Replace-RangeText $d.Paragraphs.Item(12).Range $t12 $newProgramaContent

# Avaliacao list: shift Metodo <- Criterio <- Norma <- old Bibliografia body
# Re-locate the (unchanged) bold labels fresh before each sub-write so
# earlier edits never invalidate later offsets.
$p14 = $d.Paragraphs.Item(14).Range
$mEnd = Find-LabelEnd $p14 "Método: "
$cStart = Find-LabelStart $p14 "Critério: "
Replace-RangeText ($d.Range($mEnd, $cStart)) $oldMetodoContent $oldCriterioContent

$p14 = $d.Paragraphs.Item(14).Range
$cEnd = Find-LabelEnd $p14 "Critério: "
$nStart = Find-LabelStart $p14 "Norma de recuperação: "
Replace-RangeText ($d.Range($cEnd, $nStart)) $oldCriterioContent ($oldNormaContent + $vt)

$p14 = $d.Paragraphs.Item(14).Range
$nEnd = Find-LabelEnd $p14 "Norma de recuperação: "
$pEnd = $d.Paragraphs.Item(14).Range.End
Replace-RangeText ($d.Range($nEnd, $pEnd)) $oldNormaContent $t16

# Bibliografia body <- old Docente(s) list body
Replace-RangeText $d.Paragraphs.Item(16).Range $t16 $t8
